{"js": "// Modif comparaison groupe et \u00e9quipe S\u00e9bastien\n//\n// - Remove the two red bullet paragraphs \"Sur toutes les heatmap, avoir\n//   l'option ...\" and \"Supprimer l'option de cacher les lignes du terrain\n//   sur les heatmap\".\n// - Demote the next bullet (\"Heatmap zone debut action des tirs ...\") from\n//   list level 0 to list level 1 (it keeps numId 6, just becomes a\n//   sub-bullet).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\nconst SUR_TOUTES = \"Sur toutes les \";\nconst SUPPRIMER = \"Supprimer l\\u2019option de cacher les lignes du terrain sur les \";\nconst HEATMAP_DEBUT = \"Heatmap zone debut action des tirs\";\n\nlet toDelete = [];\nlet demoteTarget = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text || \"\";\n  if (text.indexOf(SUR_TOUTES) === 0 || text.indexOf(SUPPRIMER) === 0) {\n    toDelete.push(para);\n  } else if (text.indexOf(HEATMAP_DEBUT) === 0) {\n    demoteTarget = para;\n  }\n}\n\n// Delete the two obsolete bullet paragraphs.\nfor (const para of toDelete) {\n  para.delete();\n}\n\n// Demote the following bullet to the nested list level (ilvl 1).\nif (demoteTarget) {\n  const listItem = demoteTarget.listItemOrNullObject;\n  listItem.load(\"level\");\n  await context.sync();\n  if (!listItem.isNullObject) {\n    listItem.level = 1;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Modif comparaison groupe et \u00e9quipe S\u00e9bastien\n#\n# - Remove the two red bullet paragraphs \"Sur toutes les heatmap, avoir\n#   l'option ...\" and \"Supprimer l'option de cacher les lignes du terrain\n#   sur les heatmap\".\n# - Demote the next bullet (\"Heatmap zone debut action des tirs ...\") from\n#   list level 0 to list level 1 (ListLevelNumber 1 -> 2); it keeps numId 6,\n#   just becomes a sub-bullet.\n\n$d = $word.ActiveDocument\n\n$surToutesPrefix = \"Sur toutes les \"\n$supprimerPrefix = \"Supprimer l\" + [char]0x2019 + \"option de cacher les lignes du terrain sur les \"\n$heatmapDebutPrefix = \"Heatmap zone debut action des tirs\"\n\n# Locate the paragraphs by their text (robust to any earlier edits shifting\n# paragraph indices elsewhere in the document).\n$idxSur = -1\n$idxSupprimer = -1\n$idxHeatmapDebut = -1\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.StartsWith($surToutesPrefix)) {\n        $idxSur = $i\n    } elseif ($t.StartsWith($supprimerPrefix)) {\n        $idxSupprimer = $i\n    } elseif ($t.StartsWith($heatmapDebutPrefix)) {\n        $idxHeatmapDebut = $i\n    }\n}\n\n# Delete the contiguous range spanning both obsolete bullet paragraphs in a\n# single operation (avoids stale paragraph references after a mutation).\nif ($idxSur -ne -1 -and $idxSupprimer -ne -1) {\n    $firstIdx = [Math]::Min($idxSur, $idxSupprimer)\n    $lastIdx = [Math]::Max($idxSur, $idxSupprimer)\n    $startPara = $d.Paragraphs.Item($firstIdx)\n    $endPara = $d.Paragraphs.Item($lastIdx)\n    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)\n    $delRange.Delete()\n\n    if ($idxHeatmapDebut -gt $lastIdx) {\n        $idxHeatmapDebut = $idxHeatmapDebut - 2\n    }\n}\n\n# Demote the following bullet to the nested list level (ListLevelNumber 2,\n# i.e. w:ilvl 1).\nif ($idxHeatmapDebut -ne -1) {\n    $target = $d.Paragraphs.Item($idxHeatmapDebut)\n    $target.Range.ListFormat.ListLevelNumber = 2\n}\n"}
